# Regenerate save_data to use K instead of Strike# in column G (rows 2-33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 8
    3  = 3
    4  = 10
    5  = 5
    6  = 12
    7  = 11
    8  = 5
    9  = 7
    10 = 6
    11 = 6
    12 = 5
    13 = 3
    14 = 9
    15 = 11
    16 = 8
    17 = 6
    18 = 8
    19 = 6
    20 = 3
    21 = 9
    22 = 6
    23 = 12
    24 = 9
    25 = 4
    26 = 7
    27 = 4
    28 = 8
    29 = 6
    30 = 8
    31 = 3
    32 = 5
    33 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
